$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$aValues = @(46074,46074.01041666666,46074.02083333334,46074.03125,46074.04166666666,46074.05208333334,46074.0625,46074.07291666666,46074.08333333334,46074.09375,46074.10416666666,46074.11458333334,46074.125,46074.13541666666,46074.14583333334,46074.15625,46074.16666666666,46074.17708333334,46074.1875,46074.19791666666,46074.20833333334,46074.21875,46074.22916666666,46074.23958333334,46074.25,46074.26041666666,46074.27083333334,46074.28125,46074.29166666666,46074.30208333334,46074.3125,46074.32291666666,46074.33333333334,46074.34375,46074.35416666666,46074.36458333334,46074.375,46074.38541666666,46074.39583333334,46074.40625,46074.41666666666,46074.42708333334,46074.4375,46074.44791666666,46074.45833333334,46074.46875,46074.47916666666,46074.48958333334,46074.5,46074.51041666666,46074.52083333334,46074.53125,46074.54166666666,46074.55208333334,46074.5625,46074.57291666666,46074.58333333334,46074.59375,46074.60416666666,46074.61458333334,46074.625,46074.63541666666,46074.64583333334,46074.65625,46074.66666666666,46074.67708333334,46074.6875,46074.69791666666,46074.70833333334,46074.71875,46074.72916666666,46074.73958333334,46074.75,46074.76041666666,46074.77083333334,46074.78125,46074.79166666666,46074.80208333334,46074.8125,46074.82291666666,46074.83333333334,46074.84375,46074.85416666666,46074.86458333334,46074.875,46074.88541666666,46074.89583333334,46074.90625,46074.91666666666,46074.92708333334,46074.9375,46074.94791666666,46074.95833333334,46074.96875,46074.97916666666,46074.98958333334)
$bValues = @(6230,6180,6150,6120,6090,6060,6040,6020,6010,6000,6000,6000,6010,6020,6040,6070,6110,6150,6200,6250,6310,6370,6430,6500,6580,6660,6750,6850,6950,7060,7150,7250,7330,7390,7440,7480,7500,7520,7520,7520,7510,7500,7500,7490,7480,7470,7460,7450,7440,7430,7420,7410,7410,7410,7410,7410,7410,7420,7420,7430,7440,7470,7510,7560,7620,7700,7770,7840,7900,7950,7990,8010,8020,8010,7980,7940,7880,7830,7760,7680,7580,7470,7360,7250,7110,6980,6840,6700,6580,6460,6350,6250,6230,6200,6160,6100)

$startRow = 2
for ($i = 0; $i -lt $aValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value2 = $aValues[$i]
    $ws.Cells.Item($row, 2).Value2 = $bValues[$i]
}
